# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column C ("municipio-nombre") moves from being an iaest-measure of type
# xsd:int/medida to being a curated dimension (sdmx-dimension:refArea /
# dim / URI-Municipio), mirroring column D ("provincia-nombre").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"
